$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("B2").Value = "EC302"
$wsA.Range("C2").Value = "CS251 (Elective)"
$wsA.Range("D2").Value = "MA262"

$wsA.Range("B3").Value = "CS263"
$wsA.Range("D3").Value = "CS251 (Elective)"
$wsA.Range("E3").Value = "MA262"

$wsA.Range("B5").Value = "MA261"
$wsA.Range("C5").Value = "CS263"
$wsA.Range("D5").Value = "EC301"
$wsA.Range("F5").Value = "EC301"

$wsA.Range("B6").Value = "CS251 (Tutorial)"
$wsA.Range("C6").Value = "Free"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "MA262 (Tutorial)"

$wsA.Range("C7").Value = "EC302"
$wsA.Range("D7").Value = "HS201"
$wsA.Range("E7").Value = "HS201"
$wsA.Range("F7").Value = "HS261 (Elective)"

$wsA.Range("C8").Value = "Free"
$wsA.Range("D8").Value = "EC302 (Tutorial)"
$wsA.Range("F8").Value = "CS263 (Tutorial)"

$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "CS263"
$wsB.Range("C2").Value = "CS251 (Elective)"
$wsB.Range("D2").Value = "EC301"

$wsB.Range("B3").Value = "EC302"
$wsB.Range("C3").Value = "MA262"
$wsB.Range("D3").Value = "CS251 (Elective)"
$wsB.Range("E3").Value = "MA262"

$wsB.Range("B5").Value = "MA262"
$wsB.Range("C5").Value = "EC301"
$wsB.Range("D5").Value = "HS201"
$wsB.Range("E5").Value = "EC301"

$wsB.Range("B6").Value = "CS251 (Tutorial)"
$wsB.Range("C6").Value = "MA261 (Tutorial)"
$wsB.Range("D6").Value = "MA262 (Tutorial)"

$wsB.Range("B7").Value = "MA261"
$wsB.Range("D7").Value = "CS263"
$wsB.Range("E7").Value = "EC302"
$wsB.Range("F7").Value = "HS261 (Elective)"

$wsB.Range("B8").Value = "CS263 (Tutorial)"
$wsB.Range("C8").Value = "EC302 (Tutorial)"
$wsB.Range("F8").Value = "Free"

$wsE = $wb.Worksheets.Item("Elective_Coordination")
$wsE.Range("C2").Value = "Fri"
$wsE.Range("D2").Value = "15:30-17:00"

$wsE.Range("C4").Value = "Thu"

$wsE.Range("C11").Value = "Wed"
$wsE.Range("D11").Value = "10:30-12:00"

$wsE.Range("C12").Value = "Tue"
$wsE.Range("D12").Value = "09:00-10:30"

$wsE.Range("C13").Value = "Mon"
$wsE.Range("D13").Value = "14:30-15:30"
